$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value to a range while forcing it to be stored as text,
# even when the value looks like a number (e.g. "1.002"). We temporarily
# switch the cell to Text format before assigning, then restore the
# original style so no new/visible style change is introduced.
function Set-TextValue($range, $value) {
    $style = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $style
}

# Row 46 / 47: PaxDollar and EnergySwap swapped rank position, with new
# price/volume figures.
Set-TextValue $ws.Range("B46") "EnergySwap"
Set-TextValue $ws.Range("C46") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D46") "10.17"
Set-TextValue $ws.Range("E46") "  -2.56%  "

Set-TextValue $ws.Range("B47") "PaxDollar"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue $ws.Range("D47") "1.001"
Set-TextValue $ws.Range("E47") "  +0.06%  "

Set-TextValue $ws.Range("D2") "26.443.02"
Set-TextValue $ws.Range("E2") "  -3.75%  "

Set-TextValue $ws.Range("D3") "1.770.61"
Set-TextValue $ws.Range("E3") "  -2.99%  "

Set-TextValue $ws.Range("D4") "1.002"
Set-TextValue $ws.Range("E4") "  +0.05%  "

Set-TextValue $ws.Range("D5") "1.001"
Set-TextValue $ws.Range("E5") "  +0.02%  "

Set-TextValue $ws.Range("D6") "306.52"
Set-TextValue $ws.Range("E6") "  -2.06%  "

Set-TextValue $ws.Range("D7") "0.4299"
Set-TextValue $ws.Range("E7") "  +0.93%  "

Set-TextValue $ws.Range("D8") "0.3664"
Set-TextValue $ws.Range("E8") "  +1.42%  "

Set-TextValue $ws.Range("D9") "0.07201"
Set-TextValue $ws.Range("E9") "  -0.09%  "

Set-TextValue $ws.Range("D10") "0.8500"
Set-TextValue $ws.Range("E10") "  -1.55%  "

Set-TextValue $ws.Range("E11") "  -1.27%  "

Set-TextValue $ws.Range("D12") "1.785.05"
Set-TextValue $ws.Range("E12") "  -5.38%  "

Set-TextValue $ws.Range("D13") "6.433"
Set-TextValue $ws.Range("E13") "  -0.65%  "

Set-TextValue $ws.Range("D14") "5.237"
Set-TextValue $ws.Range("E14") "  -2.95%  "

Set-TextValue $ws.Range("D15") "0.06867"
Set-TextValue $ws.Range("E15") "  -0.82%  "

Set-TextValue $ws.Range("D16") "1.005"
Set-TextValue $ws.Range("E16") "  +0.15%  "

Set-TextValue $ws.Range("D17") "79.42"
Set-TextValue $ws.Range("E17") "  -1.96%  "

Set-TextValue $ws.Range("D18") "0.000008654"
Set-TextValue $ws.Range("E18") "  -3.24%  "

Set-TextValue $ws.Range("E19") "  +0.14%  "

Set-TextValue $ws.Range("D20") "15.01"
Set-TextValue $ws.Range("E20") "  -2.29%  "

Set-TextValue $ws.Range("D21") "26.444.59"
Set-TextValue $ws.Range("E21") "  -4.63%  "

Set-TextValue $ws.Range("D22") "5.105"
Set-TextValue $ws.Range("E22") "  -0.45%  "

Set-TextValue $ws.Range("D23") "11.27"
Set-TextValue $ws.Range("E23") "  +3.69%  "

Set-TextValue $ws.Range("D24") "2.005.07"
Set-TextValue $ws.Range("E24") "  -4.08%  "

Set-TextValue $ws.Range("D25") "152.15"
Set-TextValue $ws.Range("E25") "  -2.06%  "

Set-TextValue $ws.Range("D26") "1.876"
Set-TextValue $ws.Range("E26") "  -5.71%  "

Set-TextValue $ws.Range("D27") "18.13"
Set-TextValue $ws.Range("E27") "  -3.26%  "

Set-TextValue $ws.Range("D28") "5.090"
Set-TextValue $ws.Range("E28") "  -1.16%  "

Set-TextValue $ws.Range("D29") "114.57"
Set-TextValue $ws.Range("E29") "  +0.09%  "

Set-TextValue $ws.Range("D30") "1.725"
Set-TextValue $ws.Range("E30") "  -3.63%  "

Set-TextValue $ws.Range("D31") "0.08951"
Set-TextValue $ws.Range("E31") "  +0.60%  "

Set-TextValue $ws.Range("D32") "0.7252"
Set-TextValue $ws.Range("E32") "  -3.16%  "

Set-TextValue $ws.Range("E33") "  -0.34%  "

Set-TextValue $ws.Range("D34") "4.323"
Set-TextValue $ws.Range("E34") "  -4.86%  "

Set-TextValue $ws.Range("E35") "  +0.07%  "

Set-TextValue $ws.Range("D36") "2.743"
Set-TextValue $ws.Range("E36") "  -7.71%  "

Set-TextValue $ws.Range("D37") "1.077"
Set-TextValue $ws.Range("E37") "  -0.73%  "

Set-TextValue $ws.Range("D38") "0.05166"
Set-TextValue $ws.Range("E38") "  -2.04%  "

Set-TextValue $ws.Range("D39") "0.01893"
Set-TextValue $ws.Range("E39") "  -1.86%  "

Set-TextValue $ws.Range("D40") "0.4930"

Set-TextValue $ws.Range("D41") "0.1604"
Set-TextValue $ws.Range("E41") "  -3.25%  "

Set-TextValue $ws.Range("D42") "2.543"
Set-TextValue $ws.Range("E42") "  -8.92%  "

Set-TextValue $ws.Range("D43") "6.245"
Set-TextValue $ws.Range("E43") "  -2.93%  "

Set-TextValue $ws.Range("D44") "8.028"
Set-TextValue $ws.Range("E44") "  -3.86%  "

Set-TextValue $ws.Range("D45") "104.96"
Set-TextValue $ws.Range("E45") "  -1.39%  "

Set-TextValue $ws.Range("D48") "0.06198"
Set-TextValue $ws.Range("E48") "  -4.08%  "

Set-TextValue $ws.Range("D49") "0.4482"
Set-TextValue $ws.Range("E49") "  -4.39%  "

Set-TextValue $ws.Range("E50") "  -1.82%  "

Set-TextValue $ws.Range("D51") "1.750"
Set-TextValue $ws.Range("E51") "  +2.76%  "
